$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.165.16'
$ws.Range("E2").Value = '  +0.19%  '

# Row 3
$ws.Range("D3").Value = '1.828.19'
$ws.Range("E3").Value = '  -0.57%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.43%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6194'
$ws.Range("D6").Style = "Normal"

# Row 7
$ws.Range("E7").Value = '  +0.04%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07355'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.13%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2916'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.92%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.08'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.24%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07670'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.42%  '

# Row 12
$ws.Range("D12").Value = '1.821.19'
$ws.Range("E12").Value = '  -0.56%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.948'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.36%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6637'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.85%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.16'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.90%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008906'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.92%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.832'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.27%  '

# Row 18
$ws.Range("D18").Value = '29.121.70'
$ws.Range("E18").Value = '  +0.08%  '

# Row 19
$ws.Range("D19").Value = '2.065.85'
$ws.Range("E19").Value = '  -0.46%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '238.39'
$ws.Range("D20").Style = "Normal"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.87%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.02%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.345'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.53%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.12%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.89'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.48%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1415'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.07%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.508'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.41%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.64'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.42%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.485'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.62%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05917'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.32%  '

# Row 31
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.068'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.85%  '

# Row 32
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.080'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.52%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.208'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.00%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.865'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.83%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7312'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.37%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.135'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.08%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.612'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.97%  '

# Row 38
$ws.Range("E38").Value = '  +2.74%  '

# Row 39
$ws.Range("D39").Value = '1.216.27'
$ws.Range("E39").Value = '  -1.70%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01749'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.71%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.285'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.73%  '

# Row 42
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9165'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.22%  '

# Row 43
$ws.Range("E43").Value = '  +0.09%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.87'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.58%  '

# Row 45
$ws.Range("D45").Value = '1.969.52'
$ws.Range("E45").Value = '  -0.55%  '

# Row 46
$ws.Range("E46").Value = '  -2.57%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5087'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.10%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.146'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.24%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4015'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.49%  '

# Row 50
$ws.Range("E50").Value = '  -4.58%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1129'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.83%  '
